$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1033
$ws.Range("I18").Value = 1033
$ws.Range("K18").Value = 1033
$ws.Range("M18").Value = -749

$ws.Range("H19").Value = 3496.75
$ws.Range("I19").Value = 1615.5
$ws.Range("J19").Value = 4303
$ws.Range("K19").Value = 1615.5
$ws.Range("L19").Value = 4303
$ws.Range("M19").Value = -1440.5
$ws.Range("N19").Value = -4653

$ws.Range("H33").Value = 327.53845
$ws.Range("I33").Value = 379.77777
$ws.Range("K33").Value = 379.77777
$ws.Range("M33").Value = -150.77777

$ws.Range("H41").Value = 509.64706
$ws.Range("I41").Value = 194
$ws.Range("K41").Value = 194
$ws.Range("M41").Value = 246

$ws.Range("H43").Value = 6412.091
$ws.Range("I43").Value = 4999.8423
$ws.Range("K43").Value = 4999.8423
$ws.Range("M43").Value = -4930.8423

$ws.Range("H108").Value = 120000
$ws.Range("J108").Value = 120000
$ws.Range("L108").Value = 120000
$ws.Range("N108").Value = -127680

$ws.Range("H120").Value = 190000
$ws.Range("J120").Value = 190000
$ws.Range("L120").Value = 190000
$ws.Range("N120").Value = -199676

$ws.Range("H121").Value = 1021.55554
$ws.Range("J121").Value = 1021.55554
$ws.Range("L121").Value = 3064.66662
$ws.Range("N121").Value = -6558.66662

$ws.Range("H125").Value = 2069.3572
$ws.Range("I125").Value = 766.8
$ws.Range("J125").Value = 2793
$ws.Range("K125").Value = 6901.2
$ws.Range("L125").Value = 25137
$ws.Range("M125").Value = -4441.2
$ws.Range("N125").Value = -30057

$ws.Range("H137").Value = 6003.048
$ws.Range("I137").Value = 2962.1538
$ws.Range("K137").Value = 8886.4614
$ws.Range("M137").Value = -6336.4614

$ws.Range("H138").Value = 3211.4197
$ws.Range("I138").Value = 2706.5
$ws.Range("J138").Value = 3376.9673
$ws.Range("K138").Value = 8119.5
$ws.Range("L138").Value = 10130.9019
$ws.Range("M138").Value = -2979.5
$ws.Range("N138").Value = -20410.9019

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 212666.67
$ws.Range("I34").Value = 212666.67
$ws.Range("K34").Value = 212666.67
$ws.Range("M34").Value = -212395.67

$ws.Range("H74").Value = 241295.33
$ws.Range("I74").Value = 314054.12
$ws.Range("J74").Value = 8467.200000000001
$ws.Range("K74").Value = 314054.12
$ws.Range("L74").Value = 8467.200000000001
$ws.Range("M74").Value = -313180.12
$ws.Range("N74").Value = -10215.2

$ws.Range("H77").Value = 241295.33
$ws.Range("I77").Value = 314054.12
$ws.Range("J77").Value = 8467.200000000001
$ws.Range("K77").Value = 1570270.6
$ws.Range("L77").Value = 42336
$ws.Range("M77").Value = -1565902.6
$ws.Range("N77").Value = -51072

$ws.Range("H82").Value = 38988
$ws.Range("J82").Value = 38988
$ws.Range("L82").Value = 38988
$ws.Range("N82").Value = -39710

$ws.Range("H85").Value = 38988
$ws.Range("J85").Value = 38988
$ws.Range("L85").Value = 38988
$ws.Range("N85").Value = -41484

$ws.Range("H122").Value = 3204.65
$ws.Range("I122").Value = 2878.1428
$ws.Range("J122").Value = 3966.5
$ws.Range("K122").Value = 8634.428400000001
$ws.Range("L122").Value = 11899.5
$ws.Range("M122").Value = -6184.428400000001
$ws.Range("N122").Value = -16799.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2953.325
$ws.Range("I20").Value = 2506.2
$ws.Range("K20").Value = 2506.2
$ws.Range("M20").Value = -2259.2

$ws.Range("H117").Value = 109313.664
$ws.Range("J117").Value = 109313.664
$ws.Range("L117").Value = 109313.664
$ws.Range("N117").Value = -118491.664

$ws.Range("H119").Value = 34855.5
$ws.Range("J119").Value = 34855.5
$ws.Range("L119").Value = 34855.5
$ws.Range("N119").Value = -44531.5

$ws.Range("H120").Value = 93836.664
$ws.Range("J120").Value = 93836.664
$ws.Range("L120").Value = 93836.664
$ws.Range("N120").Value = -103512.664

$ws.Range("H134").Value = 3795.889
$ws.Range("I134").Value = 2695.375
$ws.Range("K134").Value = 8086.125
$ws.Range("M134").Value = -5551.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 21033.666
$ws.Range("I36").Value = 21024
$ws.Range("J36").Value = 21053
$ws.Range("K36").Value = 21024
$ws.Range("L36").Value = 21053
$ws.Range("M36").Value = -20636
$ws.Range("N36").Value = -21829

$ws.Range("H40").Value = 21033.666
$ws.Range("I40").Value = 21024
$ws.Range("J40").Value = 21053
$ws.Range("K40").Value = 21024
$ws.Range("L40").Value = 21053
$ws.Range("M40").Value = -20864
$ws.Range("N40").Value = -21373

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9959.857
$ws.Range("I2").Value = 44.5
$ws.Range("J2").Value = 13926
$ws.Range("K2").Value = 267
$ws.Range("L2").Value = 83556
$ws.Range("M2").Value = -154
$ws.Range("N2").Value = -83782

$ws.Range("H4").Value = 7116249

$ws.Range("H11").Value = 33337750
$ws.Range("I11").Value = 5973.8423
$ws.Range("J11").Value = 90910820
$ws.Range("K11").Value = 17921.5269
$ws.Range("L11").Value = 272732460
$ws.Range("M11").Value = -17781.5269
$ws.Range("N11").Value = -272732740

$ws.Range("H55").Value = 2216.3635
$ws.Range("J55").Value = 2138
$ws.Range("L55").Value = 6414
$ws.Range("N55").Value = -6768

$ws.Range("H112").Value = 2511
$ws.Range("I112").Value = 2511
$ws.Range("K112").Value = 7533
$ws.Range("M112").Value = -6425

$ws.Range("H118").Value = 1507.25
$ws.Range("I118").Value = 1507.25
$ws.Range("K118").Value = 4521.75
$ws.Range("M118").Value = -3278.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 25928.691
$ws.Range("I2").Value = 335.66666
$ws.Range("J2").Value = 111238.78
$ws.Range("K2").Value = 335.66666
$ws.Range("L2").Value = 111238.78
$ws.Range("M2").Value = -222.66666
$ws.Range("N2").Value = -111464.78

$ws.Range("H80").Value = 5689.8096
$ws.Range("I80").Value = 5469.6665
$ws.Range("K80").Value = 5469.6665
$ws.Range("M80").Value = -4471.6665

$ws.Range("H83").Value = 5689.8096
$ws.Range("I83").Value = 5469.6665
$ws.Range("K83").Value = 27348.3325
$ws.Range("M83").Value = -22356.3325

$ws.Range("H113").Value = 3805.2104
$ws.Range("I113").Value = 3027.1667
$ws.Range("K113").Value = 3027.1667
$ws.Range("M113").Value = -857.1667000000002

$ws.Range("H122").Value = 6666.7856
$ws.Range("J122").Value = 9232.25
$ws.Range("L122").Value = 27696.75
$ws.Range("N122").Value = -32596.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3122.3489
$ws.Range("I22").Value = 2265.5
$ws.Range("J22").Value = 3867.4348
$ws.Range("K22").Value = 2265.5
$ws.Range("L22").Value = 3867.4348
$ws.Range("M22").Value = -1970.5
$ws.Range("N22").Value = -4457.4348

$ws.Range("H27").Value = 3122.3489
$ws.Range("I27").Value = 2265.5
$ws.Range("J27").Value = 3867.4348
$ws.Range("K27").Value = 2265.5
$ws.Range("L27").Value = 3867.4348
$ws.Range("M27").Value = -2158.5
$ws.Range("N27").Value = -4081.4348

$ws.Range("H46").Value = 7121.5
$ws.Range("J46").Value = 9405.4
$ws.Range("L46").Value = 9405.4
$ws.Range("N46").Value = -9781.4

$ws.Range("H55").Value = 251.97437
$ws.Range("I55").Value = 191.85715
$ws.Range("J55").Value = 405
$ws.Range("K55").Value = 191.85715
$ws.Range("L55").Value = 405
$ws.Range("M55").Value = -18.85714999999999
$ws.Range("N55").Value = -751

$ws.Range("H63").Value = 51736.5
$ws.Range("J63").Value = 53982
$ws.Range("L63").Value = 53982
$ws.Range("N63").Value = -55480

$ws.Range("H66").Value = 51736.5
$ws.Range("J66").Value = 53982
$ws.Range("L66").Value = 161946
$ws.Range("N66").Value = -169434

$ws.Range("H93").Value = 1863.5454
$ws.Range("I93").Value = 1900
$ws.Range("K93").Value = 1900
$ws.Range("M93").Value = -652

$ws.Range("H118").Value = 37995
$ws.Range("J118").Value = 37995
$ws.Range("L118").Value = 37995
$ws.Range("N118").Value = -41309

$ws.Range("H127").Value = 94999.5
$ws.Range("J127").Value = 94999.5
$ws.Range("L127").Value = 94999.5
$ws.Range("N127").Value = -104919.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 48999
$ws.Range("J16").Value = 48999
$ws.Range("L16").Value = 48999
$ws.Range("N16").Value = -49583

$ws.Range("H116").Value = 113898.75
$ws.Range("J116").Value = 113898.75
$ws.Range("L116").Value = 113898.75
$ws.Range("N116").Value = -123076.75

$ws.Range("H122").Value = 3662.9375
$ws.Range("I122").Value = 3400.5
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 10201.5
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -7751.5
$ws.Range("N122").Value = -21400

$ws.Range("H132").Value = 5560.952
$ws.Range("I132").Value = 1979.3077
$ws.Range("K132").Value = 5937.9231
$ws.Range("M132").Value = -3407.9231

$ws.Range("H136").Value = 3403485
$ws.Range("J136").Value = 3326.3076
$ws.Range("L136").Value = 9978.9228
$ws.Range("N136").Value = -15078.9228
